# CryCompanywiseStockReport_1.xlsx update:
# Quantities (col F) were revised down for a batch of stock-report line
# items, their line Values (col G = Rate(D) * Qty(F)) recomputed, the
# per-company "Sub Total:" (col B) rows and the grand totals (B619/B620)
# rolled up again, and two pairs of adjustment rows (227/228, 229/230)
# had their B/E/F/G figures swapped between the pair.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F5").Value = 22
$ws.Range("G5").Value = 8289.82
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 29.57
$ws.Range("B10").Value = 31106.79
$ws.Range("F25").Value = 46
$ws.Range("G25").Value = 1413.12
$ws.Range("B32").Value = 13244.1
$ws.Range("F64").Value = 123
$ws.Range("G64").Value = 9983.91
$ws.Range("F70").Value = 36
$ws.Range("G70").Value = 4858.2
$ws.Range("F86").Value = 81
$ws.Range("G86").Value = 10163.07
$ws.Range("B90").Value = 201482.09
$ws.Range("F115").Value = 230
$ws.Range("G115").Value = 22266.3
$ws.Range("B117").Value = 16318.58
$ws.Range("F144").Value = 1224
$ws.Range("G144").Value = 10342.8
$ws.Range("F145").Value = 656
$ws.Range("G145").Value = 5241.44
$ws.Range("B147").Value = 18615.08
$ws.Range("F151").Value = 101
$ws.Range("G151").Value = 8774.879999999999
$ws.Range("F153").Value = 32
$ws.Range("G153").Value = 1487.68
$ws.Range("B156").Value = 35636.7
$ws.Range("B227").Value = 55373
$ws.Range("E227").Value = 163.62
$ws.Range("F227").Value = -94
$ws.Range("G227").Value = -13562.32
$ws.Range("B228").Value = 63520
$ws.Range("E228").Value = 153.4
$ws.Range("F228").Value = 67
$ws.Range("G228").Value = 9666.76
$ws.Range("B229").Value = 57802
$ws.Range("E229").Value = 162.71
$ws.Range("F229").Value = -79
$ws.Range("G229").Value = -11334.92
$ws.Range("B230").Value = 63531
$ws.Range("E230").Value = 152.53
$ws.Range("F230").Value = 67
$ws.Range("G230").Value = 9613.16
$ws.Range("F247").Value = 157
$ws.Range("G247").Value = 16313.87
$ws.Range("F255").Value = 608
$ws.Range("G255").Value = 104168.64
$ws.Range("F256").Value = 293
$ws.Range("G256").Value = 44292.81
$ws.Range("B260").Value = 209722.74
$ws.Range("F280").Value = 146
$ws.Range("G280").Value = 24694.44
$ws.Range("F284").Value = 0
$ws.Range("G284").Value = 0
$ws.Range("F285").Value = 15
$ws.Range("G285").Value = 418.95
$ws.Range("F303").Value = 40
$ws.Range("G303").Value = 8435.6
$ws.Range("B304").Value = 196976.59
$ws.Range("F320").Value = 71
$ws.Range("G320").Value = 4874.15
$ws.Range("F328").Value = 59
$ws.Range("G328").Value = 2195.39
$ws.Range("B330").Value = 32092.07
$ws.Range("F339").Value = 6
$ws.Range("G339").Value = 284.4
$ws.Range("F342").Value = 142
$ws.Range("G342").Value = 4497.14
$ws.Range("F345").Value = 79
$ws.Range("G345").Value = 4851.39
$ws.Range("B346").Value = 28649.43
$ws.Range("F454").Value = 52
$ws.Range("G454").Value = 1775.8
$ws.Range("B460").Value = 15156.98
$ws.Range("F508").Value = 60
$ws.Range("G508").Value = 6236.4
$ws.Range("B510").Value = 26251.02
$ws.Range("F551").Value = 7
$ws.Range("G551").Value = 1001.91
$ws.Range("F555").Value = 35
$ws.Range("G555").Value = 2434.6
$ws.Range("B560").Value = 9076.77
$ws.Range("F599").Value = 2065
$ws.Range("G599").Value = 336822.15
$ws.Range("F602").Value = 352
$ws.Range("G602").Value = 50916.8
$ws.Range("B606").Value = 521253.03
$ws.Range("F617").Value = 2
$ws.Range("G617").Value = 79.06
$ws.Range("B618").Value = 46948.53
$ws.Range("B619").Value = 1990545.4
$ws.Range("B620").Value = 1990545.4
